$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Case_1")

# Row 53: medical cost dropped from 1,800,000 to 720,000 (cascades through P53/Q53 and everything below)
$ws.Range("O53").Value = 720000

# Row 57: category changed from "cn" to "lb"; add medical-cost-carryover column O57
$ws.Range("L57").Value = "lb"
$ws.Range("O57").Formula = "=O55"

# Row 58: new "cn" (accident) step - add category + medical cost O58
$ws.Range("L58").Value = "cn"
$ws.Range("O58").Formula = "=720000+O55"

# Row 59: new "tra" step - add category + medical cost O59
$ws.Range("L59").Value = "tra"
$ws.Range("O59").Formula = "=1800000+O55"

# Row 60: new "ca" step - add category + medical cost O60
$ws.Range("L60").Value = "ca"
$ws.Range("O60").Formula = "=O55"

# Row 61: new "die" step - add category + medical cost O61
$ws.Range("L61").Value = "die"
$ws.Range("O61").Formula = "=O60"

# Update the view to match where the author left the selection
$ws.Activate()
$ws.Range("O62").Select()
